$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.143.05"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.57%  '

$ws.Range("D3").Value = "'2.279.25"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.28%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").Value = "'154.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +15,380.32%  '

$ws.Range("D6").Value = "'305.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.34%  '

$ws.Range("D7").Value = "'94.14"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.76%  '

$ws.Range("E8").Value = '  -0.42%  '

$ws.Range("E9").Value = '  -0.02%  '

$ws.Range("E10").Value = '  +0.58%  '

$ws.Range("D11").Value = "'34.07"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.83%  '

$ws.Range("E12").Value = '  +0.21%  '

$ws.Range("E13").Value = '  -2.26%  '

$ws.Range("E14").Value = '  -0.55%  '

$ws.Range("D15").Value = "'2.632.96"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.32%  '

$ws.Range("D16").Value = "'14.37"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.49%  '

$ws.Range("D17").Value = "'2.275.56"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.80%  '

$ws.Range("E18").Value = '  +3.71%  '

$ws.Range("D19").Value = "'42.051.76"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.53%  '

$ws.Range("D20").Value = "'12.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.13%  '

$ws.Range("E21").Value = '  +1.02%  '

$ws.Range("E22").Value = '  +0.61%  '

$ws.Range("D23").Value = "'68.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.03%  '

$ws.Range("D24").Value = "'243.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.16%  '

$ws.Range("D25").Value = "'2.60"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.65%  '

$ws.Range("E26").Value = '  +0.66%  '

$ws.Range("D28").Value = "'24.11"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.89%  '

$ws.Range("E29").Value = '  +4.97%  '

$ws.Range("D30").Value = "'9.70"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.54%  '

$ws.Range("E31").Value = '  +1.36%  '

$ws.Range("D32").Value = "'160.78"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.46%  '

$ws.Range("E33").Value = '  +2.92%  '

$ws.Range("E34").Value = '  +0.01%  '

$ws.Range("D35").Value = "'0.0753"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.02%  '

$ws.Range("E36").Value = '  +0.88%  '

$ws.Range("D37").Value = "'0.109"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.56%  '

$ws.Range("D38").Value = "'17.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.00%  '

$ws.Range("E39").Value = '  -0.52%  '

$ws.Range("E40").Value = '  -0.18%  '

$ws.Range("E41").Value = '  -0.73%  '

$ws.Range("D42").Value = "'4.24"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +7.59%  '

$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").Value = "'2.022.68"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.54%  '

$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").Value = "'19.76"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.35%  '

$ws.Range("E45").Value = '  +11.41%  '

$ws.Range("E46").Value = '  +1.34%  '

$ws.Range("E47").Value = '  -2.15%  '

$ws.Range("E48").Value = '  +0.28%  '

$ws.Range("E49").Value = '  +3.23%  '

$ws.Range("E50").Value = '  -1.30%  '

$ws.Range("D51").Value = "'72.19"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.46%  '
